$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the shared header labels: "_old" -> "_FV2410" and "_new" -> "_FV2504"
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value
    if ($v -ne $null) {
        if ($v.ToString().EndsWith("_old")) {
            $cell.Value = $v.ToString().Replace("_old", "_FV2410")
        } elseif ($v.ToString().EndsWith("_new")) {
            $cell.Value = $v.ToString().Replace("_new", "_FV2504")
        }
    }
}

# Turn the data range into an Excel Table ("Table1")
$tableRange = $ws.Range("A1:U76")
$list = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$list.Name = "Table1"

# Freeze the header row
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
